# "Added a gimic to call the student's name when they login"
#
# The student_pswd sheet used column A purely as a throw-away row index
# (1..9) that isn't needed any more now that a name look-up exists (see
# column C / the G5:G6 "num rows" helper). Clear those stale index
# numbers, bump the helper count in G6, and leave the selection where
# the author last clicked (G9) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("student_pswd")
$ws.Activate()

# Drop the old sequential row-index values in column A (rows 4-10).
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()

# "num rows" helper (G5 label / G6 value) now reflects 2 rows.
$ws.Range("G6").Value = 2

# Leave the cursor on G9, matching where the workbook was saved from.
$ws.Range("G9").Select()

# Restore the workbook window's last on-screen position/size.
$win = $excel.ActiveWindow
$win.Left = 17325
$win.Top = 5295
$win.Width = 28800
$win.Height = 16335
